$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated (shortened) description text used across all product rows.
$desc = "Erkek baggy pantolon, bol ve rahat kesimiyle öne çıkan, modern sokak stilinin vazgeçilmez parçasıdır. Kalçadan paçaya kadar geniş formu sayesinde hareket özgürlüğü sunar.Ürün içeriği 100% Pamuk.29-38  beden seçeneği mevcuttur."

# Refresh the "aciklama" column on the existing rows with the shortened text.
$ws.Cells.Item(2, 5).Value = $desc
$ws.Cells.Item(3, 5).Value = $desc
$ws.Cells.Item(4, 5).Value = $desc
$ws.Cells.Item(5, 5).Value = $desc

# Fix the price casing on the "Baggy Füme" row (300 tl -> 300 Tl).
$ws.Cells.Item(5, 2).Value = "300 Tl"

# Add the new "Baggy Kar Yıkama" product row.
$ws.Cells.Item(6, 1).Value = "Baggy Kar Yıkama"
$ws.Cells.Item(6, 2).Value = "300 Tl"
$ws.Cells.Item(6, 3).Value = "Jeans"
$ws.Cells.Item(6, 4).Value = "BAG4.jpg"
$ws.Cells.Item(6, 5).Value = $desc
$ws.Cells.Item(6, 6).Value = "Var"

# Widen column E slightly to fit the new content.
$ws.Columns.Item(5).ColumnWidth = 9.7

# Move the active selection to F6, matching the latest edit location.
$ws.Range("F6").Select() | Out-Null
